# Update NATMI Arf1-Insr LR-pairs sheet with refreshed TPM-derived statistics.
# The underlying TPM recomputation changed the per-(sending-cluster) ligand
# average/total expression (columns G/H) and the per-(target-cluster) receptor
# average/total expression (columns M/N); every other touched column
# (I,J,O,P,Q,R,S,T) is a value derived from those via the NATMI specificity /
# edge-weight formulas (I=G/sum(G), J=H/sum(H), O=M/sum(M), P=N/sum(N),
# Q=G*M, R=H*N, S=I*O, T=J*P). The new literal values below were taken from
# the target workbook's canonical XML and are written directly, matching the
# original file's storage of plain values (no formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 49.88947433333334
$ws.Cells.Item(2, 8).Value2 = 149.668423
$ws.Cells.Item(2, 9).Value2 = 0.2324880572195875
$ws.Cells.Item(2, 10).Value2 = 0.2324880572195874
$ws.Cells.Item(2, 13).Value2 = 20.29529466666667
$ws.Cells.Item(2, 14).Value2 = 60.885884
$ws.Cells.Item(2, 15).Value2 = 0.4032332285476398
$ws.Cells.Item(2, 16).Value2 = 0.4032332285476398
$ws.Cells.Item(2, 17).Value2 = 1012.521582360104
$ws.Cells.Item(2, 18).Value2 = 9112.694241240933
$ws.Cells.Item(2, 19).Value2 = 0.09374690991142266
$ws.Cells.Item(2, 20).Value2 = 0.09374690991142265
$ws.Cells.Item(3, 7).Value2 = 49.88947433333334
$ws.Cells.Item(3, 8).Value2 = 149.668423
$ws.Cells.Item(3, 9).Value2 = 0.2324880572195875
$ws.Cells.Item(3, 10).Value2 = 0.2324880572195874
$ws.Cells.Item(3, 15).Value2 = 0.1953894087318433
$ws.Cells.Item(3, 16).Value2 = 0.1953894087318433
$ws.Cells.Item(3, 17).Value2 = 490.6242325766012
$ws.Cells.Item(3, 18).Value2 = 4415.618093189411
$ws.Cells.Item(3, 19).Value2 = 0.04542570403735014
$ws.Cells.Item(3, 20).Value2 = 0.04542570403735013
$ws.Cells.Item(4, 7).Value2 = 49.88947433333334
$ws.Cells.Item(4, 8).Value2 = 149.668423
$ws.Cells.Item(4, 9).Value2 = 0.2324880572195875
$ws.Cells.Item(4, 10).Value2 = 0.2324880572195874
$ws.Cells.Item(4, 13).Value2 = 11.81535133333333
$ws.Cells.Item(4, 14).Value2 = 35.446054
$ws.Cells.Item(4, 15).Value2 = 0.2347510761885954
$ws.Cells.Item(4, 16).Value2 = 0.2347510761885954
$ws.Cells.Item(4, 17).Value2 = 589.4616670836491
$ws.Cells.Item(4, 18).Value2 = 5305.155003752842
$ws.Cells.Item(4, 19).Value2 = 0.05457682163329389
$ws.Cells.Item(4, 20).Value2 = 0.05457682163329389
$ws.Cells.Item(5, 7).Value2 = 49.88947433333334
$ws.Cells.Item(5, 8).Value2 = 149.668423
$ws.Cells.Item(5, 9).Value2 = 0.2324880572195875
$ws.Cells.Item(5, 10).Value2 = 0.2324880572195874
$ws.Cells.Item(5, 13).Value2 = 8.386535
$ws.Cells.Item(5, 14).Value2 = 25.159605
$ws.Cells.Item(5, 15).Value2 = 0.1666262865319216
$ws.Cells.Item(5, 16).Value2 = 0.1666262865319216
$ws.Cells.Item(5, 17).Value2 = 418.3998226281017
$ws.Cells.Item(5, 18).Value2 = 3765.598403652915
$ws.Cells.Item(5, 19).Value2 = 0.03873862163752077
$ws.Cells.Item(5, 20).Value2 = 0.03873862163752075
$ws.Cells.Item(6, 9).Value2 = 0.295249080025651
$ws.Cells.Item(6, 10).Value2 = 0.295249080025651
$ws.Cells.Item(6, 13).Value2 = 20.29529466666667
$ws.Cells.Item(6, 14).Value2 = 60.885884
$ws.Cells.Item(6, 15).Value2 = 0.4032332285476398
$ws.Cells.Item(6, 16).Value2 = 0.4032332285476398
$ws.Cells.Item(6, 17).Value2 = 1285.855580166767
$ws.Cells.Item(6, 18).Value2 = 11572.7002215009
$ws.Cells.Item(6, 19).Value2 = 0.1190542397644637
$ws.Cells.Item(6, 20).Value2 = 0.1190542397644637
$ws.Cells.Item(7, 9).Value2 = 0.295249080025651
$ws.Cells.Item(7, 10).Value2 = 0.295249080025651
$ws.Cells.Item(7, 15).Value2 = 0.1953894087318433
$ws.Cells.Item(7, 16).Value2 = 0.1953894087318433
$ws.Cells.Item(7, 19).Value2 = 0.05768854317483263
$ws.Cells.Item(7, 20).Value2 = 0.05768854317483263
$ws.Cells.Item(8, 9).Value2 = 0.295249080025651
$ws.Cells.Item(8, 10).Value2 = 0.295249080025651
$ws.Cells.Item(8, 13).Value2 = 11.81535133333333
$ws.Cells.Item(8, 14).Value2 = 35.446054
$ws.Cells.Item(8, 15).Value2 = 0.2347510761885954
$ws.Cells.Item(8, 16).Value2 = 0.2347510761885954
$ws.Cells.Item(8, 17).Value2 = 748.5890544151833
$ws.Cells.Item(8, 18).Value2 = 6737.301489736649
$ws.Cells.Item(8, 19).Value2 = 0.0693100392797143
$ws.Cells.Item(8, 20).Value2 = 0.0693100392797143
$ws.Cells.Item(9, 9).Value2 = 0.295249080025651
$ws.Cells.Item(9, 10).Value2 = 0.295249080025651
$ws.Cells.Item(9, 13).Value2 = 8.386535
$ws.Cells.Item(9, 14).Value2 = 25.159605
$ws.Cells.Item(9, 15).Value2 = 0.1666262865319216
$ws.Cells.Item(9, 16).Value2 = 0.1666262865319216
$ws.Cells.Item(9, 17).Value2 = 531.348423618875
$ws.Cells.Item(9, 18).Value2 = 4782.135812569875
$ws.Cells.Item(9, 19).Value2 = 0.04919625780664038
$ws.Cells.Item(9, 20).Value2 = 0.04919625780664037
$ws.Cells.Item(10, 7).Value2 = 52.37451933333333
$ws.Cells.Item(10, 8).Value2 = 157.123558
$ws.Cells.Item(10, 9).Value2 = 0.2440685216737345
$ws.Cells.Item(10, 10).Value2 = 0.2440685216737345
$ws.Cells.Item(10, 13).Value2 = 20.29529466666667
$ws.Cells.Item(10, 14).Value2 = 60.885884
$ws.Cells.Item(10, 15).Value2 = 0.4032332285476398
$ws.Cells.Item(10, 16).Value2 = 0.4032332285476398
$ws.Cells.Item(10, 17).Value2 = 1062.95630289503
$ws.Cells.Item(10, 18).Value2 = 9566.606726055274
$ws.Cells.Item(10, 19).Value2 = 0.09841653798134956
$ws.Cells.Item(10, 20).Value2 = 0.09841653798134956
$ws.Cells.Item(11, 7).Value2 = 52.37451933333333
$ws.Cells.Item(11, 8).Value2 = 157.123558
$ws.Cells.Item(11, 9).Value2 = 0.2440685216737345
$ws.Cells.Item(11, 10).Value2 = 0.2440685216737345
$ws.Cells.Item(11, 15).Value2 = 0.1953894087318433
$ws.Cells.Item(11, 16).Value2 = 0.1953894087318433
$ws.Cells.Item(11, 17).Value2 = 515.0627200999845
$ws.Cells.Item(11, 18).Value2 = 4635.56448089986
$ws.Cells.Item(11, 19).Value2 = 0.04768840413988606
$ws.Cells.Item(11, 20).Value2 = 0.04768840413988606
$ws.Cells.Item(12, 7).Value2 = 52.37451933333333
$ws.Cells.Item(12, 8).Value2 = 157.123558
$ws.Cells.Item(12, 9).Value2 = 0.2440685216737345
$ws.Cells.Item(12, 10).Value2 = 0.2440685216737345
$ws.Cells.Item(12, 13).Value2 = 11.81535133333333
$ws.Cells.Item(12, 14).Value2 = 35.446054
$ws.Cells.Item(12, 15).Value2 = 0.2347510761885954
$ws.Cells.Item(12, 16).Value2 = 0.2347510761885954
$ws.Cells.Item(12, 17).Value2 = 618.8233468377923
$ws.Cells.Item(12, 18).Value2 = 5569.410121540132
$ws.Cells.Item(12, 19).Value2 = 0.05729534812666869
$ws.Cells.Item(12, 20).Value2 = 0.05729534812666869
$ws.Cells.Item(13, 7).Value2 = 52.37451933333333
$ws.Cells.Item(13, 8).Value2 = 157.123558
$ws.Cells.Item(13, 9).Value2 = 0.2440685216737345
$ws.Cells.Item(13, 10).Value2 = 0.2440685216737345
$ws.Cells.Item(13, 13).Value2 = 8.386535
$ws.Cells.Item(13, 14).Value2 = 25.159605
$ws.Cells.Item(13, 15).Value2 = 0.1666262865319216
$ws.Cells.Item(13, 16).Value2 = 0.1666262865319216
$ws.Cells.Item(13, 17).Value2 = 439.2407394971767
$ws.Cells.Item(13, 18).Value2 = 3953.16665547459
$ws.Cells.Item(13, 19).Value2 = 0.0406682314258302
$ws.Cells.Item(13, 20).Value2 = 0.0406682314258302
$ws.Cells.Item(14, 7).Value2 = 48.96808833333333
$ws.Cells.Item(14, 8).Value2 = 146.904265
$ws.Cells.Item(14, 9).Value2 = 0.2281943410810271
$ws.Cells.Item(14, 10).Value2 = 0.228194341081027
$ws.Cells.Item(14, 13).Value2 = 20.29529466666667
$ws.Cells.Item(14, 14).Value2 = 60.885884
$ws.Cells.Item(14, 15).Value2 = 0.4032332285476398
$ws.Cells.Item(14, 16).Value2 = 0.4032332285476398
$ws.Cells.Item(14, 17).Value2 = 993.8217819883622
$ws.Cells.Item(14, 18).Value2 = 8944.396037895262
$ws.Cells.Item(14, 19).Value2 = 0.09201554089040385
$ws.Cells.Item(14, 20).Value2 = 0.09201554089040384
$ws.Cells.Item(15, 7).Value2 = 48.96808833333333
$ws.Cells.Item(15, 8).Value2 = 146.904265
$ws.Cells.Item(15, 9).Value2 = 0.2281943410810271
$ws.Cells.Item(15, 10).Value2 = 0.228194341081027
$ws.Cells.Item(15, 15).Value2 = 0.1953894087318433
$ws.Cells.Item(15, 16).Value2 = 0.1953894087318433
$ws.Cells.Item(15, 17).Value2 = 481.5631168763945
$ws.Cells.Item(15, 18).Value2 = 4334.068051887551
$ws.Cells.Item(15, 19).Value2 = 0.04458675737977445
$ws.Cells.Item(15, 20).Value2 = 0.04458675737977445
$ws.Cells.Item(16, 7).Value2 = 48.96808833333333
$ws.Cells.Item(16, 8).Value2 = 146.904265
$ws.Cells.Item(16, 9).Value2 = 0.2281943410810271
$ws.Cells.Item(16, 10).Value2 = 0.228194341081027
$ws.Cells.Item(16, 13).Value2 = 11.81535133333333
$ws.Cells.Item(16, 14).Value2 = 35.446054
$ws.Cells.Item(16, 15).Value2 = 0.2347510761885954
$ws.Cells.Item(16, 16).Value2 = 0.2347510761885954
$ws.Cells.Item(16, 17).Value2 = 578.5751677800343
$ws.Cells.Item(16, 18).Value2 = 5207.17651002031
$ws.Cells.Item(16, 19).Value2 = 0.0535688671489185
$ws.Cells.Item(16, 20).Value2 = 0.05356886714891849
$ws.Cells.Item(17, 7).Value2 = 48.96808833333333
$ws.Cells.Item(17, 8).Value2 = 146.904265
$ws.Cells.Item(17, 9).Value2 = 0.2281943410810271
$ws.Cells.Item(17, 10).Value2 = 0.228194341081027
$ws.Cells.Item(17, 13).Value2 = 8.386535
$ws.Cells.Item(17, 14).Value2 = 25.159605
$ws.Cells.Item(17, 15).Value2 = 0.1666262865319216
$ws.Cells.Item(17, 16).Value2 = 0.1666262865319216
$ws.Cells.Item(17, 17).Value2 = 410.6725866905917
$ws.Cells.Item(17, 18).Value2 = 3696.053280215325
$ws.Cells.Item(17, 19).Value2 = 0.03802317566193027
$ws.Cells.Item(17, 20).Value2 = 0.03802317566193025
